# User_Schema_Rights_Definition.xlsx edit
# Commit: "DB - Creating views to display the latest version of the FHIR resources"
#
# Inserts 4 new rows (40-43) on "rights_and_functions" describing new
# "_last_version" view-creation scripts (typed/raw x cds2db/dataproc),
# mirroring the existing "_last" view rows. Existing rows 40-55 shift to 44-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert 4 blank rows above the existing row 40 ------------------
$ws.Rows.Item(40).Resize(4).Insert()

# --- 2. Populate the new rows --------------------------------------------
# Set cells in a specific order so newly-created shared strings are
# interned in the same order as the target workbook.
$ws.Range("C40").Value = "template_cre_view5.sql"
$ws.Range("H41").Value = "_raw_last_version"
$ws.Range("H40").Value = "_last_version"
$ws.Range("B40").Value = "23_cre_view_typ_cds2db_last_version.sql"
$ws.Range("B41").Value = "23_cre_view_raw_cds2db_last_version.sql"
$ws.Range("B42").Value = "23_cre_view_typ_dataproc_last_version.sql"
$ws.Range("B43").Value = "23_cre_view_raw_dataproc_last_version.sql"

# Row 40: 23_cre_view_typ_cds2db_last_version.sql
$ws.Range("D40").Value = "cds2db_in"
$ws.Range("E40").Value = "Spaltenname"
$ws.Range("G40").Value = "v_"
$ws.Range("I40").Value = "SELECT"
$ws.Range("J40").Value = "cds2db_in"
$ws.Range("N40").Value = "db_log"

# Row 41: 23_cre_view_raw_cds2db_last_version.sql
$ws.Range("C41").Value = "template_cre_view5.sql"
$ws.Range("D41").Value = "cds2db_in"
$ws.Range("E41").Value = "Spaltenname"
$ws.Range("G41").Value = "v_"
$ws.Range("I41").Value = "SELECT"
$ws.Range("J41").Value = "cds2db_in"
$ws.Range("N41").Value = "db_log"
$ws.Range("O41").Value = "_raw"

# Row 42: 23_cre_view_typ_dataproc_last_version.sql
$ws.Range("C42").Value = "template_cre_view5.sql"
$ws.Range("D42").Value = "db2dataprocessor_user"
$ws.Range("E42").Value = "db2dataprocessor_out"
$ws.Range("G42").Value = "v_"
$ws.Range("H42").Value = "_last_version"
$ws.Range("I42").Value = "SELECT"
$ws.Range("J42").Value = "db2dataprocessor_user"
$ws.Range("N42").Value = "db_log"

# Row 43: 23_cre_view_raw_dataproc_last_version.sql
$ws.Range("C43").Value = "template_cre_view5.sql"
$ws.Range("D43").Value = "db2dataprocessor_user"
$ws.Range("E43").Value = "db2dataprocessor_out"
$ws.Range("G43").Value = "v_"
$ws.Range("H43").Value = "_raw_last_version"
$ws.Range("I43").Value = "SELECT"
$ws.Range("J43").Value = "db2dataprocessor_user"
$ws.Range("N43").Value = "db_log"
$ws.Range("O43").Value = "_raw"

# --- 3. Re-home the two cell comments that were pushed down by the insert ---
$commentText = $ws.Range("K46").Comment.Text()
$ws.Range("K46").Comment.Delete()
$ws.Range("K50").AddComment($commentText)

$commentText2 = $ws.Range("K48").Comment.Text()
$ws.Range("K48").Comment.Delete()
$ws.Range("K52").AddComment($commentText2)

# --- 4. Update the view/selection state ----------------------------------
$ws.Range("E42").Select()
